$wb = $excel.ActiveWorkbook

# --- 1. Insert a new "2022-Q1" sheet, positioned right after "2021-Q4" ---
$sheetQ4 = $wb.Worksheets.Item("2021-Q4")

$newSheet = $wb.Worksheets.Add($null, $sheetQ4)
$newSheet.Name = "2022-Q1"

# Re-fetch handles by name after the structural change above: worksheet
# references captured before an Add()/Insert() can end up pointing at the
# wrong physical sheet once the collection is re-shuffled.
$sheetQ4 = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Item("2022-Q1")

# Copy the layout/styling of the 2021-Q4 sheet (header row + 2 data rows)
# as the starting point for the new 2022-Q1 sheet.
$sheetQ4.Range("A1:H3").Copy($newSheet.Range("A1"))
$newSheet.Range("A1").ClearContents()

# Overwrite with the 2022-Q1 fund-holding figures. These columns hold
# numeric-looking text (matching the source data's inline-string typing),
# so force text entry via NumberFormat, then restore the default "Normal"
# style afterwards so no stray text-format style lingers on the cells.
$newSheet.Range("D2:G3").NumberFormat = "@"

$newSheet.Range("D2").Value = "16.02"
$newSheet.Range("E2").Value = "93.69"
$newSheet.Range("F2").Value = "4.00"
$newSheet.Range("G2").Value = "0.6408"
$newSheet.Range("H2").Value = 8

$newSheet.Range("D3").Value = "0.51"
$newSheet.Range("E3").Value = "93.69"
$newSheet.Range("F3").Value = "4.00"
$newSheet.Range("G3").Value = "0.0204"
$newSheet.Range("H3").Value = 8

$newSheet.Range("D2:G3").Style = "Normal"

# --- 2. Update the "总计" (summary) sheet: add a 2022-Q1 row on top, ---
#        pushing the existing 2021-Q4 row down.
$sheetTotal = $wb.Worksheets.Item("总计")
$sheetTotal.Rows.Item(2).Insert()
$sheetTotal = $wb.Worksheets.Item("总计")

# Re-use row 3's style (carried forward from the original row 2) for the
# freshly-inserted row 2, then fix up the index values.
$sheetTotal.Range("A3").Copy($sheetTotal.Range("A2"))
$sheetTotal.Range("A2").Value = 0
$sheetTotal.Range("A3").Value = 1

$sheetTotal.Range("B2:D2").ClearFormats()
$sheetTotal.Range("B2").Value = "2022-Q1"
$sheetTotal.Range("C2").Value = 2
$sheetTotal.Range("D2").Value = 0.66
